# Automatic sync of changes
$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- Create "Sheet2" (holds the grad-program tracker table) ---
$sheet2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet2.Name = "Sheet2"

# --- Create "Sheet3" (stays empty, ends up as the active tab) ---
$sheet3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet3.Name = "Sheet3"

# Move Sheet3 so the tab order becomes Sheet1, Sheet3, Sheet2
$sheet3.Move($sheet2)

# Re-fetch handles by name now that positions have shifted
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet3 = $wb.Worksheets.Item("Sheet3")

# --- Populate "Sheet2" with the grad-program tracker data ---
$headers = @("id", "university", "degree", "program", "mode", "sent email", "application deadline", "funding", "stipend", "semester start date")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $sheet2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$sheet2.Cells.Item(2, 1).Value = 1
$sheet2.Cells.Item(2, 2).Value = "Harvard University"
$sheet2.Cells.Item(2, 3).Value = "PhD"
$sheet2.Cells.Item(2, 4).Value = "biomedical informatics"
$sheet2.Cells.Item(2, 5).Value = "on-campus"
$sheet2.Cells.Item(2, 6).Value = "yes"

$sheet2.Cells.Item(3, 1).Value = 2
$sheet2.Cells.Item(3, 2).Value = "Columbia University"
$sheet2.Cells.Item(3, 3).Value = "PhD"
$sheet2.Cells.Item(3, 4).Value = "biomedical informatics"
$sheet2.Cells.Item(3, 5).Value = "on-campus"
$sheet2.Cells.Item(3, 6).Value = "yes"

$sheet2.Cells.Item(4, 1).Value = 3
$sheet2.Cells.Item(4, 2).Value = "University of Utah"
$sheet2.Cells.Item(4, 3).Value = "PhD"
$sheet2.Cells.Item(4, 4).Value = "biomedical informatics"
$sheet2.Cells.Item(4, 5).Value = "hybrid"
$sheet2.Cells.Item(4, 6).Value = "yes"

$sheet2.Cells.Item(5, 1).Value = 4
$sheet2.Cells.Item(5, 2).Value = "George Mason University"
$sheet2.Cells.Item(5, 3).Value = "PhD"
$sheet2.Cells.Item(5, 4).Value = "bioinformatics"
$sheet2.Cells.Item(5, 5).Value = "hybrid"
$sheet2.Cells.Item(5, 6).Value = "yes"
$sheet2.Cells.Item(5, 7).Value = 45962
$sheet2.Cells.Item(5, 7).NumberFormat = "mm-dd-yy"
$sheet2.Cells.Item(5, 8).Value = "full"

$sheet2.Cells.Item(6, 1).Value = 5
$sheet2.Cells.Item(6, 2).Value = "Stanford University"
$sheet2.Cells.Item(6, 3).Value = "PhD"
$sheet2.Cells.Item(6, 4).Value = "biomedical informatics"
$sheet2.Cells.Item(6, 5).Value = "on-campus"
$sheet2.Cells.Item(6, 6).Value = "yes - waiting for reply"

$sheet2.Cells.Item(7, 2).Value = "Purdue University"
$sheet2.Cells.Item(7, 3).Value = "PhD"

$sheet2.Cells.Item(8, 2).Value = "Purdue University"
$sheet2.Cells.Item(8, 3).Value = "MSc"
$sheet2.Cells.Item(8, 4).Value = "health informatics"

$sheet2.Cells.Item(7, 4).Value = "health sciences"

# Column widths on "Sheet2" (best achievable approximation of the
# Excel "best fit" auto-sized widths; this engine quantizes column widths
# to 1/6-character pixel units, so the closest representable value is used)
$sheet2.Columns.Item(1).ColumnWidth = 1.5
$sheet2.Columns.Item(2).ColumnWidth = 21.166666666666668
$sheet2.Columns.Item(3).ColumnWidth = 5.666666666666667
$sheet2.Columns.Item(4).ColumnWidth = 19
$sheet2.Columns.Item(5).ColumnWidth = 9.333333333333334
$sheet2.Columns.Item(6).ColumnWidth = 18
$sheet2.Columns.Item(7).ColumnWidth = 17
$sheet2.Columns.Item(8).ColumnWidth = 6.333333333333333
$sheet2.Columns.Item(9).ColumnWidth = 6.166666666666667
$sheet2.Columns.Item(10).ColumnWidth = 16.333333333333332

$sheet2.Range("D7").Select()

# --- Sheet1 tweaks: clear the bold style that used to live on B1/A2 ---
$sheet1.Range("B1").Font.Bold = $false
$sheet1.Range("A2").Font.Bold = $false

$sheet1.Range("B17").Select()

# --- "Sheet3" (empty) becomes the active tab ---
$sheet3.Activate()
$sheet3.Range("A1").Select()
